# chore: update Sheets via scheduled runner
#
# Refreshes the market-price-derived columns (H:N) for the affected leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets, mirroring the latest
# pull from the scheduled price-scraper run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2416.5
$ws.Range("I20").Value = 1899.8
$ws.Range("J20").Value = 5000
$ws.Range("K20").Value = 1899.8
$ws.Range("L20").Value = 5000
$ws.Range("M20").Value = -1669.8
$ws.Range("N20").Value = -5460

$ws.Range("H33").Value = 50185.535
$ws.Range("I33").Value = 62677.832
$ws.Range("K33").Value = 62677.832
$ws.Range("M33").Value = -62448.832

$ws.Range("H35").Value = 2416.5
$ws.Range("I35").Value = 1899.8
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 1899.8
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -1520.8
$ws.Range("N35").Value = -5758

$ws.Range("H69").Value = 6043.364
$ws.Range("J69").Value = 6043.364
$ws.Range("L69").Value = 18130.092
$ws.Range("N69").Value = -19878.092

$ws.Range("H72").Value = 6043.364
$ws.Range("J72").Value = 6043.364
$ws.Range("L72").Value = 54390.276
$ws.Range("N72").Value = -63126.276

$ws.Range("H103").Value = 731.0769
$ws.Range("I103").Value = 529.4286
$ws.Range("K103").Value = 1588.2858
$ws.Range("M103").Value = -1002.2858

$ws.Range("H132").Value = 4843.766
$ws.Range("I132").Value = 4523.5815
$ws.Range("J132").Value = 8285.75
$ws.Range("K132").Value = 13570.7445
$ws.Range("L132").Value = 24857.25
$ws.Range("M132").Value = -11040.7445
$ws.Range("N132").Value = -29917.25

$ws.Range("H138").Value = 6039.88
$ws.Range("I138").Value = 5019.2144
$ws.Range("J138").Value = 6436.8057
$ws.Range("K138").Value = 15057.6432
$ws.Range("L138").Value = 19310.4171
$ws.Range("M138").Value = -9917.643199999999
$ws.Range("N138").Value = -29590.4171

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21488
$ws.Range("I32").Value = 17758.076
$ws.Range("J32").Value = 45732.5
$ws.Range("K32").Value = 17758.076
$ws.Range("L32").Value = 45732.5
$ws.Range("M32").Value = -17471.076
$ws.Range("N32").Value = -46306.5

$ws.Range("H74").Value = 2668.5122
$ws.Range("I74").Value = 2114.8708
$ws.Range("J74").Value = 4384.8
$ws.Range("K74").Value = 2114.8708
$ws.Range("L74").Value = 4384.8
$ws.Range("M74").Value = -1240.8708
$ws.Range("N74").Value = -6132.8

$ws.Range("H77").Value = 2668.5122
$ws.Range("I77").Value = 2114.8708
$ws.Range("J77").Value = 4384.8
$ws.Range("K77").Value = 10574.354
$ws.Range("L77").Value = 21924
$ws.Range("M77").Value = -6206.354000000001
$ws.Range("N77").Value = -30660

$ws.Range("H94").Value = 30000
$ws.Range("J94").Value = 30000
$ws.Range("L94").Value = 30000
$ws.Range("N94").Value = -31802

$ws.Range("H135").Value = 65166
$ws.Range("J135").Value = 65166
$ws.Range("L135").Value = 65166
$ws.Range("N135").Value = -75306

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4282.7144
$ws.Range("I134").Value = 3955.4119
$ws.Range("J134").Value = 5673.75
$ws.Range("K134").Value = 11866.2357
$ws.Range("L134").Value = 17021.25
$ws.Range("M134").Value = -9331.235700000001
$ws.Range("N134").Value = -22091.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 4799.3335
$ws.Range("I3").Value = 3999
$ws.Range("J3").Value = 5199.5
$ws.Range("K3").Value = 3999
$ws.Range("L3").Value = 5199.5
$ws.Range("M3").Value = -3886
$ws.Range("N3").Value = -5425.5

$ws.Range("H31").Value = 4726.304
$ws.Range("I31").Value = 4105.8184
$ws.Range("J31").Value = 5295.0835
$ws.Range("K31").Value = 4105.8184
$ws.Range("L31").Value = 5295.0835
$ws.Range("M31").Value = -3810.8184
$ws.Range("N31").Value = -5885.0835

$ws.Range("H34").Value = 4726.304
$ws.Range("I34").Value = 4105.8184
$ws.Range("J34").Value = 5295.0835
$ws.Range("K34").Value = 4105.8184
$ws.Range("L34").Value = 5295.0835
$ws.Range("M34").Value = -3903.8184
$ws.Range("N34").Value = -5699.0835

$ws.Range("H45").Value = 13516.5
$ws.Range("I45").Value = 13022
$ws.Range("K45").Value = 13022
$ws.Range("M45").Value = -12429

$ws.Range("H58").Value = 361252.78
$ws.Range("I58").Value = 2506.3
$ws.Range("J58").Value = 560556.4
$ws.Range("K58").Value = 2506.3
$ws.Range("L58").Value = 560556.4
$ws.Range("M58").Value = -2303.3
$ws.Range("N58").Value = -560962.4

$ws.Range("H136").Value = 361252.78
$ws.Range("I136").Value = 2506.3
$ws.Range("J136").Value = 560556.4
$ws.Range("K136").Value = 7518.900000000001
$ws.Range("L136").Value = 1681669.2
$ws.Range("M136").Value = -4968.900000000001
$ws.Range("N136").Value = -1686769.2

$ws.Range("H141").Value = 264024.75
$ws.Range("J141").Value = 271115.7
$ws.Range("L141").Value = 271115.7
$ws.Range("N141").Value = -281475.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 773.6
$ws.Range("I2").Value = 57.5
$ws.Range("K2").Value = 345
$ws.Range("M2").Value = -232

$ws.Range("H5").Value = 669.3125
$ws.Range("I5").Value = 617.0833
$ws.Range("K5").Value = 1851.2499
$ws.Range("M5").Value = -1739.2499

$ws.Range("H15").Value = 21.75
$ws.Range("I15").Value = 25.8
$ws.Range("K15").Value = 77.40000000000001
$ws.Range("M15").Value = 62.59999999999999

$ws.Range("H16").Value = 288
$ws.Range("I16").Value = 240
$ws.Range("J16").Value = 297.6
$ws.Range("K16").Value = 720
$ws.Range("L16").Value = 892.8000000000001
$ws.Range("M16").Value = -547
$ws.Range("N16").Value = -1238.8

$ws.Range("H26").Value = 251.66667
$ws.Range("I26").Value = 225
$ws.Range("K26").Value = 675
$ws.Range("M26").Value = -387

$ws.Range("H32").Value = 2387.5
$ws.Range("J32").Value = 3000
$ws.Range("L32").Value = 9000
$ws.Range("N32").Value = -9566

$ws.Range("H122").Value = 1074.0769
$ws.Range("I122").Value = 909
$ws.Range("J122").Value = 1445.5
$ws.Range("K122").Value = 8181
$ws.Range("L122").Value = 13009.5
$ws.Range("M122").Value = -5731
$ws.Range("N122").Value = -17909.5

$ws.Range("H131").Value = 19544.924
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()

$ws.Range("H132").Value = 1857.9333
$ws.Range("I132").Value = 1583.4286
$ws.Range("J132").Value = 2098.125
$ws.Range("K132").Value = 14250.8574
$ws.Range("L132").Value = 18883.125
$ws.Range("M132").Value = -11720.8574
$ws.Range("N132").Value = -23943.125

$ws.Range("H135").Value = 669.3125
$ws.Range("I135").Value = 617.0833
$ws.Range("K135").Value = 5553.7497
$ws.Range("M135").Value = -3018.7497

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5083333.5
$ws.Range("I7").Value = 2625000
$ws.Range("J7").Value = 10000000
$ws.Range("K7").Value = 2625000
$ws.Range("L7").Value = 10000000
$ws.Range("M7").Value = -2624888
$ws.Range("N7").Value = -10000224

$ws.Range("H8").Value = 5083333.5
$ws.Range("I8").Value = 2625000
$ws.Range("J8").Value = 10000000
$ws.Range("K8").Value = 2625000
$ws.Range("L8").Value = 10000000
$ws.Range("M8").Value = -2624861
$ws.Range("N8").Value = -10000278

$ws.Range("H28").Value = 13999
$ws.Range("I28").Value = 13999
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 13999
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -13807
$ws.Range("N28").ClearContents()

$ws.Range("H126").Value = 7914.727
$ws.Range("I126").Value = 14209.667
$ws.Range("J126").Value = 5554.125
$ws.Range("K126").Value = 42629.001
$ws.Range("L126").Value = 16662.375
$ws.Range("M126").Value = -40159.001
$ws.Range("N126").Value = -21602.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 45458930
$ws.Range("I7").Value = 83336450
$ws.Range("K7").Value = 83336450
$ws.Range("M7").Value = -83336338

$ws.Range("H122").Value = 9032.25
$ws.Range("I122").Value = 8980.333000000001
$ws.Range("K122").Value = 26940.999
$ws.Range("M122").Value = -24490.999

$ws.Range("H126").Value = 45458930
$ws.Range("I126").Value = 83336450
$ws.Range("K126").Value = 250009350
$ws.Range("M126").Value = -250006880

$ws.Range("H132").Value = 120942.07
$ws.Range("I132").Value = 185172.45
$ws.Range("J132").Value = 6984.9355
$ws.Range("K132").Value = 555517.3500000001
$ws.Range("L132").Value = 20954.8065
$ws.Range("M132").Value = -552987.3500000001
$ws.Range("N132").Value = -26014.8065

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 13339272
$ws.Range("I81").Value = 4958.7
$ws.Range("J81").Value = 40007900
$ws.Range("K81").Value = 9917.4
$ws.Range("L81").Value = 80015800
$ws.Range("M81").Value = -8856.4
$ws.Range("N81").Value = -80017922

$ws.Range("H84").Value = 13339272
$ws.Range("I84").Value = 4958.7
$ws.Range("J84").Value = 40007900
$ws.Range("K84").Value = 49587
$ws.Range("L84").Value = 400079000
$ws.Range("M84").Value = -44283
$ws.Range("N84").Value = -400089608

$ws.Range("H132").Value = 290098.5
$ws.Range("I132").Value = 349243.7
$ws.Range("J132").Value = 4230
$ws.Range("K132").Value = 1047731.1
$ws.Range("L132").Value = 12690
$ws.Range("M132").Value = -1045201.1
$ws.Range("N132").Value = -17750
